# "Added new taskings for GUI"
# Adds two new backlog rows (GUI tasking items) to the "Sprint 1" backlog
# sheet, extends the trailing blank rows down to row 29, refreshes the
# "Total" formula, and leaves the "Sprint 1" tab as the active/selected
# sheet with E25 selected (previously "Burndown Chart" was the active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

# Row 17's "By" column changes from "AE" to "SS".
$ws.Range("E17").Value = "SS"

# New backlog item #15 (row 18): GUI update referencing the utility class.
$ws.Range("C18").Value = "Updated GUI to refereence utility class"
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = "AE"

# New backlog item #16 (row 19): GUI update for correct calculations.
$ws.Range("C19").Value = "Updated GUI for correct calculations"
$ws.Range("D19").Value = 2

# E19 is a brand-new cell on this row, so pick up the same formatting
# (centered, bordered style) used by the rest of column E before setting
# its value.
$ws.Range("E18").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E19").Value = "AE"

# Extend the formatted-but-empty column E cells down through the new
# bottom of the sheet (row 29), matching the existing style used at E18/E20.
$ws.Range("E18").Copy()
$ws.Range("E20:E29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# D22's Total formula (=SUM(D4:D21)) recalculates automatically to include
# the new Est hours.

# "Sprint 1" becomes the active sheet/tab, with E25 the selected cell
# (previously "Burndown Chart" was active with C18 selected).
$ws.Activate() | Out-Null
$ws.Range("E25").Select() | Out-Null
